# Update the worksheet date and every division-problem answer cell.
$d = $word.ActiveDocument

# --- Update the date/day heading paragraph ---
$d.Content.Find.Execute("2024-08-22 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-23 Friday", 2)

# --- Update each answer cell in the (only) table by explicit row/column ---
# Using Cell.Range.Text avoids ambiguity from duplicate "old" values
# (e.g. "60÷9=6, 6" and "75÷8=9, 3" each occur twice with different
# replacements) that a document-wide Find/Replace could not disambiguate.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "72÷4=18, 0"
$t.Cell(1, 2).Range.Text  = "51÷4=12, 3"
$t.Cell(1, 3).Range.Text  = "99÷7=14, 1"
$t.Cell(1, 4).Range.Text  = "11÷3=3, 2"
$t.Cell(1, 5).Range.Text  = "51÷8=6, 3"

$t.Cell(5, 1).Range.Text  = "54÷9=6, 0"
$t.Cell(5, 2).Range.Text  = "27÷4=6, 3"
$t.Cell(5, 3).Range.Text  = "60÷4=15, 0"
$t.Cell(5, 4).Range.Text  = "56÷9=6, 2"
$t.Cell(5, 5).Range.Text  = "10÷6=1, 4"

$t.Cell(9, 1).Range.Text  = "44÷7=6, 2"
$t.Cell(9, 2).Range.Text  = "18÷2=9, 0"
$t.Cell(9, 3).Range.Text  = "56÷5=11, 1"
$t.Cell(9, 4).Range.Text  = "67÷5=13, 2"
$t.Cell(9, 5).Range.Text  = "28÷8=3, 4"

$t.Cell(13, 1).Range.Text = "52÷6=8, 4"
$t.Cell(13, 2).Range.Text = "34÷8=4, 2"
$t.Cell(13, 3).Range.Text = "19÷4=4, 3"
$t.Cell(13, 4).Range.Text = "15÷9=1, 6"
$t.Cell(13, 5).Range.Text = "92÷9=10, 2"

$t.Cell(17, 1).Range.Text = "33÷9=3, 6"
$t.Cell(17, 2).Range.Text = "52÷9=5, 7"
$t.Cell(17, 3).Range.Text = "92÷5=18, 2"
$t.Cell(17, 4).Range.Text = "44÷6=7, 2"
$t.Cell(17, 5).Range.Text = "29÷4=7, 1"
